# Apply the updates described by the diff:
#  - Move the sheet selection from A2 to A10
#  - Update several date / percentage values in rows 9-12 and 18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("A9").Value = (Get-Date -Year 2021 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0)

# Row 10
$ws.Range("A10").Value = (Get-Date -Year 2021 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)

# Row 11
$ws.Range("A11").Value = (Get-Date -Year 2021 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D11").Value = 60
$ws.Range("E11").Value = (Get-Date -Year 2021 -Month 11 -Day 1 -Hour 0 -Minute 0 -Second 0)

# Row 12
$ws.Range("A12").Value = (Get-Date -Year 2021 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D12").Value = 70

# Row 18
$ws.Range("D18").Value = 50
$ws.Range("E18").Value = (Get-Date -Year 2021 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0)

# Move the active selection to A10 (matches the <selection activeCell="A10" sqref="A10"/> change)
$ws.Range("A10").Select()
